$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Step 1: rename the existing sheet to the new per-axle naming convention.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Tire2x_270_70R22"

# ---------------------------------------------------------------------------
# Step 2: burn a sheetId (a temporary sheet) so that the *copy* we make below
# receives the same internal sheetId ("5") that the real edit ended up with,
# then get rid of the temporary sheet again.
# ---------------------------------------------------------------------------
$tmp = $wb.Worksheets.Add()

# ---------------------------------------------------------------------------
# Step 3: duplicate the sheet to support the second axle / tire configuration.
# ---------------------------------------------------------------------------
$ws1.Copy([System.Reflection.Missing]::Value, $ws1)
$ws2 = $wb.Worksheets.Item("Tire2x_270_70R22 (2)")
$ws2.Name = "Tire2x_430_50R38"

# remove the scratch sheet used to reserve the sheetId
$wb.Worksheets.Item("Sheet1").Delete()

# put the new sheet right after the renamed original sheet
$ws2.Move([System.Reflection.Missing]::Value, $ws1)

# ---------------------------------------------------------------------------
# Step 4: update the "Tire" label cell and the instance-name cell on sheet 2
# first, so shared-string allocation order matches (Tire, Tire2x_430_50R38,
# Tire2x_270_70R22).
# ---------------------------------------------------------------------------
$ws2.Range("H2").Value = "Tire"
$ws2.Range("H3").Value = "Tire2x_430_50R38"

# ---------------------------------------------------------------------------
# Step 5: update the "Tire" label cell and the instance-name cell on sheet 1.
# ---------------------------------------------------------------------------
$ws1.Range("H2").Value = "Tire"
$ws1.Range("H3").Value = "Tire2x_270_70R22"

# ---------------------------------------------------------------------------
# Step 6: on the new (second) sheet, replace the computed xOffset formula
# with a static value representative of the new tire's specification.
# ---------------------------------------------------------------------------
$ws2.Range("H7").Value = 0.4572

# ---------------------------------------------------------------------------
# Step 7: refresh the conditional-formatting rules on sheet 1 (the renamed
# original) so that it gets its own distinct set of highlight styles,
# separate from the copy on sheet 2.
# ---------------------------------------------------------------------------
$rA20 = $ws1.Range("A20")
$rA20.FormatConditions.Item(1).Delete()
$fcA20 = $rA20.FormatConditions.Add(8, [System.Reflection.Missing]::Value, '"class"')
$fcA20.Interior.Color = 13431551

$rA19 = $ws1.Range("A19")
$rA19.FormatConditions.Item(1).Delete()
$fcA19 = $rA19.FormatConditions.Add(8, [System.Reflection.Missing]::Value, '"class"')
$fcA19.Interior.Color = 13431551

# ---------------------------------------------------------------------------
# Step 8: restore the selection / active-cell state seen in the edited file.
# ---------------------------------------------------------------------------
$ws1.Range("C25").Select()
$ws2.Range("J16").Select()

# sheet 2 was left as the active tab
$ws2.Activate()
